# Applies the "Added task to ID-1" edit:
#  - Fix wording in E3 ("and employee" -> "or employee")
#  - Add a new Task entry in F3 (the new bootstrap/AngularJS task description)
#  - Widen column F to fit the new text
#  - Grow row 3's height to fit the wrapped task text
#  - Move the active selection to F4

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Fix the small wording change in the existing user story cell.
$ws.Range("E3").Value = "As a user, I want to login to the system. Here, user can be line manager, site manager or employee."

# Add the new Task text for ID 1.
$ws.Range("F3").Value = "Is to create a webpage using bootstrap and AngularJS containing two input fields (UserName and PassWord), two buttons (Register and LogIn), a Remember me check box"

# New column width for column F to accommodate the task text.
$ws.Columns.Item(6).ColumnWidth = 29.5

# Row 3 grows to fit the new wrapped content.
$ws.Rows.Item(3).RowHeight = 106

# Update the active cell/selection.
$ws.Range("F4").Select()
